$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new command row first
$ws.Range("B13").Value = "git reset --hard origin/<nume_branch>"
$ws.Range("C13").Value = "reseteaza branch-ul local, ca sa fie ca pe online"

# Explicitly clear fill (matches the "No Fill" style applied to the new row)
$ws.Range("B13:C13").Interior.ColorIndex = -4142

# Update placeholder text on existing rows to use angle-bracket notation
$ws.Range("B4").Value = "git branch <nume_branch>"
$ws.Range("B5").Value = "git checkout -b <nume_branch>"

# Adjust the selection like the author left it
$ws.Range("B16").Select()
